$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the slip image name values (trashed items / payments renamed slips)
$ws.Range("G2").Value = "slip-image-name-1.png"
$ws.Range("G3").Value = "slip-image-name.-2png"

# Update the active selection on the sheet
$ws.Range("M5").Select()
